$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3256227758007118
$ws1.Range("C2").Value = 0.0687960687960688
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.128735632183908
$ws1.Range("F2").Value = 0.2697495183044316
$ws1.Range("G2").Value = 0.6576332429990966
$ws1.Range("H2").Value = 0.731474050294275
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 379
$ws1.Range("K2").Value = 155
$ws1.Range("L2").Value = 0

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2902621722846442
$ws2.Range("D2").Value = 0.4499274310595066

$ws2.Range("B3").Value = 0.0687960687960688
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.128735632183908

$ws2.Range("B4").Value = 0.3256227758007118
$ws2.Range("C4").Value = 0.3256227758007118
$ws2.Range("D4").Value = 0.3256227758007118
$ws2.Range("E4").Value = 0.3256227758007118

$ws2.Range("B5").Value = 0.5343980343980343
$ws2.Range("C5").Value = 0.6451310861423221
$ws2.Range("D5").Value = 0.2893315316217073

$ws2.Range("B6").Value = 0.9536054980894838
$ws2.Range("C6").Value = 0.3256227758007118
$ws2.Range("D6").Value = 0.4339249926813629

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 155
$ws3.Range("C2").Value = 379

$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
